$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- J12: copy existing "red fill" format (style 24) from G18, then set value ---
$ws.Range("G18").Copy() | Out-Null
$ws.Range("J12").PasteSpecial(-4122) | Out-Null
$ws.Range("J12").Value = "2분 24초"

# --- J13: style unchanged (s="8"), just set value ---
$ws.Range("J13").Value = "2분 47초"

# --- J14: copy style 24 (red fill) again ---
$ws.Range("G18").Copy() | Out-Null
$ws.Range("J14").PasteSpecial(-4122) | Out-Null
$ws.Range("J14").Value = "3분 10초"

# --- J15: copy style 24 (red fill) ---
$ws.Range("G18").Copy() | Out-Null
$ws.Range("J15").PasteSpecial(-4122) | Out-Null
$ws.Range("J15").Value = "3분 7초"

# --- J16: copy style 24 (red fill) ---
$ws.Range("G18").Copy() | Out-Null
$ws.Range("J16").PasteSpecial(-4122) | Out-Null
$ws.Range("J16").Value = "3분 8초"

# --- J17: copy style 24 (red fill) ---
$ws.Range("G18").Copy() | Out-Null
$ws.Range("J17").PasteSpecial(-4122) | Out-Null
$ws.Range("J17").Value = "3분 17초"

# --- J18: copy style 24 (red fill) ---
$ws.Range("G18").Copy() | Out-Null
$ws.Range("J18").PasteSpecial(-4122) | Out-Null
$ws.Range("J18").Value = "3분 22초"

# --- J19: copy style 26 (blue fill + bottom border) from F19 ---
$ws.Range("F19").Copy() | Out-Null
$ws.Range("J19").PasteSpecial(-4122) | Out-Null
$ws.Range("J19").Value = "3분 9초"

$excel.CutCopyMode = 0

# --- New row 20: short-passage practice summary line ---
$ws.Rows.Item(20).RowHeight = 13.5

$ws.Range("B20").Value = "짧은글 연습 "
$ws.Range("J20").Value = "302타 98%"

$ws.Range("B20").Font.Size = 10
$ws.Range("B20").HorizontalAlignment = -4108

$ws.Range("J20").Font.Size = 10
$ws.Range("J20").HorizontalAlignment = -4108

# --- Move active selection to the newly added cell ---
$ws.Range("J20").Select() | Out-Null
